{"js": "// Apply the two textual changes described by the diff to the first\n// body paragraph of the document:\n//   1. \"pou\u017e\u00edvan\u00fd FAT32\" -> \"pou\u017e\u00edvan\u00e9 FAT32\"\n//   2. Insert a new sentence right before the final\n//      \"Na z\u00e1ver budem prezentova\u0165 v\u00fdsledky experimentov.\" sentence:\n//      \"V neposlednom rade uk\u00e1\u017eem,c \u017ee po spr\u00e1vnom form\u00e1tovan\u00ed disku sa\n//       nebud\u00fa da\u0165 d\u00e1ta obnovi\u0165. \"\n\nconst body = context.document.body;\n\n// --- Edit 1: \"pou\u017e\u00edvan\u00fd\" -> \"pou\u017e\u00edvan\u00e9\" (agreement fix before FAT32) ---\nconst fixResults = body.search(\"pou\u017e\u00edvan\u00fd FAT32\", { matchCase: true, matchWholeWord: false });\nfixResults.load(\"text\");\nawait context.sync();\n\nif (fixResults.items.length > 0) {\n  fixResults.items[0].insertText(\"pou\u017e\u00edvan\u00e9 FAT32\", \"Replace\");\n  await context.sync();\n}\n\n// --- Edit 2: insert a new sentence before the closing sentence of the\n// first paragraph. This exact sentence text only occurs once in the\n// document (a similarly-worded \"Na z\u00e1ver ...\" sentence appears later,\n// but with different wording, so this search is unambiguous). ---\nconst closingResults = body.search(\"Na z\u00e1ver budem prezentova\u0165 v\u00fdsledky experimentov.\", { matchCase: true, matchWholeWord: false });\nclosingResults.load(\"text\");\nawait context.sync();\n\nif (closingResults.items.length > 0) {\n  closingResults.items[0].insertText(\n    \"V neposlednom rade uk\u00e1\u017eem,c \u017ee po spr\u00e1vnom form\u00e1tovan\u00ed disku sa nebud\u00fa da\u0165 d\u00e1ta obnovi\u0165. \",\n    \"Before\"\n  );\n  await context.sync();\n}\n", "ps1": "# Apply the two textual changes described by the diff to the first\n# body paragraph of the document:\n#   1. \"pou\u017e\u00edvan\u00fd FAT32\" -> \"pou\u017e\u00edvan\u00e9 FAT32\"\n#   2. Insert a new sentence right before the final\n#      \"Na z\u00e1ver budem prezentova\u0165 v\u00fdsledky experimentov.\" sentence:\n#      \"V neposlednom rade uk\u00e1\u017eem,c \u017ee po spr\u00e1vnom form\u00e1tovan\u00ed disku sa\n#       nebud\u00fa da\u0165 d\u00e1ta obnovi\u0165. \"\n\n$d = $word.ActiveDocument\n\n# --- Edit 1: \"pou\u017e\u00edvan\u00fd\" -> \"pou\u017e\u00edvan\u00e9\" (agreement fix before FAT32) ---\n$find1 = $d.Content\n$find1.Find.Execute(\"pou\u017e\u00edvan\u00fd FAT32\", $false, $false, $false, $false, $false, $true, 1, $false, \"pou\u017e\u00edvan\u00e9 FAT32\", 2)\n\n# --- Edit 2: insert a new sentence before the closing sentence of the\n# first paragraph. This exact sentence text only occurs once in the\n# document (a similarly-worded \"Na z\u00e1ver ...\" sentence appears later,\n# but with different wording, so this search is unambiguous). ---\n$find2 = $d.Content\n$ok = $find2.Find.Execute(\"Na z\u00e1ver budem prezentova\u0165 v\u00fdsledky experimentov.\")\nif ($ok) {\n    $insertRng = $find2.Duplicate\n    $insertRng.Collapse(1)\n    $insertRng.InsertBefore(\"V neposlednom rade uk\u00e1\u017eem,c \u017ee po spr\u00e1vnom form\u00e1tovan\u00ed disku sa nebud\u00fa da\u0165 d\u00e1ta obnovi\u0165. \")\n}\n"}
